# Update the "想去人数" (interested-count) figures in column F across the
# 展览 (sheet 1), 本地生活 (sheet 3), and 全部类型 (sheet 4) worksheets to
# reflect the latest counts (演出, sheet 2, is unchanged).

$wb = $excel.ActiveWorkbook

# ---- 展览 (Worksheets.Item(1)) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 252
$ws.Range("F3").Value = 2634
$ws.Range("F5").Value = 922
$ws.Range("F7").Value = 1950
$ws.Range("F8").Value = 1778
$ws.Range("F10").Value = 61
$ws.Range("F11").Value = 2434
$ws.Range("F16").Value = 114
$ws.Range("F18").Value = 8963
$ws.Range("F20").Value = 6959
$ws.Range("F21").Value = 11365
$ws.Range("F25").Value = 312
$ws.Range("F26").Value = 539
$ws.Range("F27").Value = 2488
$ws.Range("F29").Value = 187
$ws.Range("F30").Value = 2397
$ws.Range("F31").Value = 596
$ws.Range("F33").Value = 4481
$ws.Range("F34").Value = 765
$ws.Range("F35").Value = 329
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 488

# ---- 本地生活 (Worksheets.Item(3)) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 621
$ws.Range("F5").Value = 133

# ---- 全部类型 (Worksheets.Item(4)) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 621
$ws.Range("F4").Value = 252
$ws.Range("F6").Value = 2634
$ws.Range("F8").Value = 922
$ws.Range("F10").Value = 1950
$ws.Range("F12").Value = 1778
$ws.Range("F15").Value = 2434
$ws.Range("F21").Value = 114
$ws.Range("F23").Value = 8963
$ws.Range("F25").Value = 6959
$ws.Range("F26").Value = 11365
$ws.Range("F30").Value = 312
$ws.Range("F32").Value = 539
$ws.Range("F34").Value = 2488
$ws.Range("F38").Value = 187
$ws.Range("F40").Value = 4481
$ws.Range("F46").Value = 488
